{"js": "// 1. Remove the old \"_GoBack\" bookmark (currently sitting right after the\n//    \"<LD_CUC>\" placeholder, before the trailing tab).\nconst existing = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nexisting.load(\"isNullObject\");\nawait context.sync();\nif (!existing.isNullObject) {\n    context.document.deleteBookmark(\"_GoBack\");\n    await context.sync();\n}\n\n// 2. Locate the word \"t\u1ea1i\" and drop a fresh, collapsed \"_GoBack\" bookmark\n//    immediately after it (i.e. right before the run that used to read\n//    \": \" and now reads \" \").\nconst body = context.document.body;\nconst taiResults = body.search(\"t\u1ea1i\", { matchCase: false, matchWholeWord: false });\ntaiResults.load(\"items\");\nawait context.sync();\n\nconst taiRange = taiResults.items[0];\nconst afterTai = taiRange.getRange(Word.RangeLocation.end);\nafterTai.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3. Drop the colon in \"t\u1ea1i: <ten_dv>\" so it reads \"t\u1ea1i <ten_dv>\" \u2014 only\n//    the run right after \"t\u1ea1i\" changes, from \": \" to \" \".\nconst colonResults = body.search(\": <ten_dv>\", { matchCase: false, matchWholeWord: false });\ncolonResults.load(\"items\");\nawait context.sync();\n\nconst colonAndName = colonResults.items[0];\nconst parts = colonAndName.split([\":\"], false, false);\nparts.load(\"items\");\nawait context.sync();\n\nconst colonOnly = parts.items[0];\ncolonOnly.insertText(\" \", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate \"t\u1ea1i\" and the trailing \": \" run that precedes \"<ten_dv>\" without\n# relying on hard-coded character offsets.\n$taiRng = $d.Content\n$taiFind = $taiRng.Find\n$taiFind.Text = \"t\u1ea1i\"\n$taiFind.Execute() | Out-Null\n$afterTai = $taiRng.End\n\n$colonRng = $d.Content\n$colonFind = $colonRng.Find\n$colonFind.Text = \": <ten_dv>\"\n$colonFind.Execute() | Out-Null\n$colonStart = $colonRng.Start\n$colonEnd = $colonStart + 2   # just the \": \" part, not \"<ten_dv>\"\n\n# 1. Remove the old \"_GoBack\" bookmark (currently sitting right after the\n#    \"<LD_CUC>\" placeholder, before the trailing tab).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Re-insert \"_GoBack\" as an empty bookmark right after the word \"t\u1ea1i\".\n#    This also acts as a natural boundary so the upcoming text edit does not\n#    get merged into the \"t\u1ea1i\" run.\n$bmRange = $d.Range($afterTai, $afterTai)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange) | Out-Null\n\n# 3. Drop a temporary boundary bookmark right after the \": \" run so the\n#    edit below cannot ripple into (and merge with) the runs that follow\n#    (e.g. \"; M\u00e3 s\u1ed1 thu\u1ebf: \", \"<mst>\", ...).\n$d.Bookmarks.Add(\"ZZZ_TmpBoundary\", $d.Range($colonEnd, $colonEnd)) | Out-Null\n\n# 4. Drop the colon in \"t\u1ea1i: <ten_dv>\" so it reads \"t\u1ea1i <ten_dv>\" (the run\n#    right after \"t\u1ea1i\" changes from \": \" to \" \").\n$colonOnly = $d.Range($colonStart, $colonEnd)\n$colonOnly.Text = \" \"\n\n# 5. Remove the temporary boundary bookmark again.\n$d.Bookmarks(\"ZZZ_TmpBoundary\").Delete()\n"}
